$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data rows 2-14 with new optimization results ---
$ws.Range("D2").Value = 30
$ws.Range("G2").Value = 551584
$ws.Range("H2").Value = 45
$ws.Range("I2").Value = 45
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = -43058

$ws.Range("D3").Value = 35
$ws.Range("G3").Value = 468454
$ws.Range("H3").Value = 38
$ws.Range("I3").Value = 38
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = -37572

$ws.Range("D4").Value = 40
$ws.Range("G4").Value = 428688
$ws.Range("H4").Value = 35
$ws.Range("I4").Value = 35
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = -33064

$ws.Range("D5").Value = 45
$ws.Range("G5").Value = 389346
$ws.Range("H5").Value = 32
$ws.Range("I5").Value = 32
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -30176

$ws.Range("D6").Value = 50
$ws.Range("G6").Value = 354464
$ws.Range("H6").Value = 29
$ws.Range("I6").Value = 29
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = -28984

$ws.Range("D7").Value = 55
$ws.Range("G7").Value = 335906
$ws.Range("H7").Value = 27
$ws.Range("I7").Value = 27
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = -25609

$ws.Range("D8").Value = 60
$ws.Range("G8").Value = 315355
$ws.Range("H8").Value = 26
$ws.Range("I8").Value = 26
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = -25764

$ws.Range("D9").Value = 65
$ws.Range("G9").Value = 297563
$ws.Range("H9").Value = 24
$ws.Range("I9").Value = 24
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = -22792

$ws.Range("D10").Value = 70
$ws.Range("G10").Value = 267986
$ws.Range("H10").Value = 22
$ws.Range("I10").Value = 22
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = -20182

$ws.Range("D11").Value = 75
$ws.Range("G11").Value = 263962
$ws.Range("H11").Value = 22
$ws.Range("I11").Value = 22
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = -21470

$ws.Range("D12").Value = 80
$ws.Range("G12").Value = 251603
$ws.Range("H12").Value = 20
$ws.Range("I12").Value = 20
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = -20396

$ws.Range("D13").Value = 85
$ws.Range("G13").Value = 229364
$ws.Range("H13").Value = 19
$ws.Range("I13").Value = 19
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = -17128

$ws.Range("D14").Value = 90
$ws.Range("G14").Value = 218200
$ws.Range("H14").Value = 18
$ws.Range("I14").Value = 18
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = -16220

# --- Remove now-obsolete rows 15-20 ---
$ws.Rows("15:20").Delete()

# --- Apply currency number format to the Final_equity header (matches body cells) ---
$ws.Range("G1").NumberFormat = '#,##0.00\ "€"'

# --- Update selection to match final cursor position ---
$ws.Range("G14").Select()
